$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("TZP", "Year", 0.5590485170101787),
    @("TZP", "Specimen_type", 0.0007520966514480875),
    @("TZP", "Gender", 0.9523860337458645),
    @("TZP", "Age_cat", 0.4120654748312703),
    @("TZP", "Hospital:Ward_ED_ICU", 0.01069603380249411)
)

$row = 52
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
